$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Address" column header; shift old "District" header to G2
$ws.Range("G2").Value = "District"
$ws.Range("F2").Value = "Address"

# Touch H42 (without giving it a value) so the sheet's used range / dimension
# extends to column H, matching the source workbook's recorded dimension.
$ws.Range("H42").Font.Bold = $false

# Row 3
$ws.Range("G3").Value = "Mandya"
$ws.Range("F3").Value = "G H S R B HalliMalavalli"

# Row 4
$ws.Range("G4").Value = "Mandya"
$ws.Range("F4").Value = "G H S T S chathra Pandavapura"

# Row 5
$ws.Range("G5").Value = "Mandya"
$ws.Range("F5").Value = "G H S BheemanahalliNagamangala"

# Row 6
$ws.Range("G6").Value = "Mandya"
$ws.Range("F6").Value = "G H S (Hariyalamma Temple) ChittanahallyPandavapura"

# Row 7
$ws.Range("G7").Value = "Mandya"
$ws.Range("F7").Value = "DhanaguruMalavalli"

# Row 8
$ws.Range("G8").Value = "Mandy"
$ws.Range("F8").Value = "G H SchoolMuthegere Mandya North"

# Row 9
$ws.Range("G9").Value = "Mandya"
$ws.Range("F9").Value = "R E H SHanakere"

# Row 10
$ws.Range("G10").Value = "Mandya"
$ws.Range("F10").Value = "G G H S (New) Malavalli"

# Row 11
$ws.Range("G11").Value = "Mandya"
$ws.Range("F11").Value = "Govt. High SchoolB G PuraMalavally"

# Row 12
$ws.Range("G12").Value = "Mandya"
$ws.Range("F12").Value = "Abhinava BharathiVidya Kerndra High School"

# Row 13
$ws.Range("G13").Value = "Mandya"
$ws.Range("F13").Value = "St. Joseph High School M C Road"

# Row 14
$ws.Range("G14").Value = "Mandya – 57145"
$ws.Range("F14").Value = "Shri Shambhulingeshwara High School Yere GowdanahallyPandavapura"

# Row 15
$ws.Range("G15").Value = "Mandya"
$ws.Range("F15").Value = "B R H S RamenahalliNagamanagala"

# Row 16
$ws.Range("G16").Value = "Mandya"
$ws.Range("F16").Value = "G H S HeraganahalliNagamangala(Tq)"

# Row 17
$ws.Range("G17").Value = "Mandya"
$ws.Range("F17").Value = "G H S MudalakoppaluPandavapura"

# Row 18
$ws.Range("G18").Value = "Mandya"
$ws.Range("F18").Value = "G H S NaguvanahalliSri Rangapatna"

# Row 19
$ws.Range("G19").Value = "Mandya"
$ws.Range("F19").Value = "G H S AlambadiK R Pete"

# Row 20
$ws.Range("G20").Value = "Mandya"
$ws.Range("F20").Value = "G P U C Dudda"

# Row 21
$ws.Range("G21").Value = "Mandya"
$ws.Range("F21").Value = "Bharathi High SchoolK M DoddiMaddur"

# Row 22
$ws.Range("G22").Value = "Mandya"
$ws.Range("F22").Value = "G P U CK R Pet"

# Row 23
$ws.Range("G23").Value = "SALEEM PASHA "
$ws.Range("F23").ClearContents()

# Row 24
$ws.Range("G24").Value = "Mandya"
$ws.Range("F24").Value = "Jnanasurya High SchoolBelakavadiMalavalli"

# Row 25
$ws.Range("G25").Value = "Mandya"
$ws.Range("F25").Value = "G H S YathambadiMalavalli"

# Row 26
$ws.Range("G26").Value = "G H S Makavalli K R Pet(Tq) Mandya(Dist."
$ws.Range("F26").ClearContents()

# Row 27
$ws.Range("G27").Value = "Mandya"
$ws.Range("F27").Value = "G H S A Hullukere"

# Row 28
$ws.Range("G28").Value = "Mandya"
$ws.Range("F28").Value = "G J C Chinakurali Pandavapura"

# Row 29
$ws.Range("G29").Value = "G H S Uramarakasalagere Mandya"
$ws.Range("F29").ClearContents()

# Row 30
$ws.Range("G30").Value = "Mandya"
$ws.Range("F30").Value = "G H S G Malligere"

# Row 31
$ws.Range("G31").Value = "Mandya"
$ws.Range("F31").Value = "Sri Yadushyla High School MelukotePandavapura"

# Row 32
$ws.Range("G32").Value = "Mandya"
$ws.Range("F32").Value = "G H S DoddakothagereMandya North"

# Row 33
$ws.Range("G33").Value = "Mandya"
$ws.Range("F33").Value = "G H S BidarahosallyMaddur"

# Row 34
$ws.Range("G34").Value = "Mandya"
$ws.Range("F34").Value = "Sarvajanika High School BelakavadiMalavalli"

# Row 35
$ws.Range("G35").Value = "Mandya"
$ws.Range("F35").Value = "Alk G C BellurNagamangala"

# Row 36
$ws.Range("G36").Value = "Mandya"
$ws.Range("F36").Value = "G J C Bindiganavile Nagamangala"

# Row 37
$ws.Range("G37").Value = "Mandya"
$ws.Range("F37").Value = "G P U C KodiyalaS R Patna"

# Row 38
$ws.Range("G38").Value = "Mandya"
$ws.Range("F38").Value = "Govt. High School K ShettallySrirangapatna"

# Row 39
$ws.Range("G39").Value = "Mandya North"
$ws.Range("F39").Value = "G H S Hulikere"

# Row 40
$ws.Range("G40").Value = "Mandya"
$ws.Range("F40").Value = "Govt. Junior college High School Section ChikkanakanahalliMaddur"

# Row 41
$ws.Range("G41").Value = "Mandya"
$ws.Range("F41").Value = "Govt. High School GananguruS R Patna Tq"

# Row 42
$ws.Range("G42").Value = "Mandya"
$ws.Range("F42").Value = "G H S KatteriPandavapura"
